$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round row 5 values down to 2 decimal places ("custom accuracy")
$ws.Range("B5").Value = 11.33
$ws.Range("C5").Value = 8.85
$ws.Range("D5").Value = 0.55
$ws.Range("E5").Value = 24.74
$ws.Range("F5").Value = 20.6
$ws.Range("G5").Value = 9.09
$ws.Range("H5").Value = 30.7
$ws.Range("I5").Value = 13.91
$ws.Range("J5").Value = 6.23
$ws.Range("K5").Value = 9.73
$ws.Range("L5").Value = 10.73
$ws.Range("M5").Value = 10.5
$ws.Range("N5").Value = 2.73
$ws.Range("O5").Value = 8.79
$ws.Range("P5").Value = 12.86
$ws.Range("Q5").Value = 7.45
$ws.Range("R5").Value = 0.08
$ws.Range("S5").Value = 0.31
$ws.Range("T5").Value = 129.8
$ws.Range("U5").Value = 25
$ws.Range("V5").Value = 8.04
$ws.Range("W5").Value = 16.77
$ws.Range("X5").Value = 9.23
$ws.Range("Y5").Value = 1.2
$ws.Range("Z5").Value = 15.72
$ws.Range("AA5").Value = 7.26
$ws.Range("AB5").Value = 6.96
$ws.Range("AC5").Value = 7.41
$ws.Range("AD5").Value = 10.93
$ws.Range("AE5").Value = 0.08
$ws.Range("AF5").Value = 27.42
$ws.Range("AG5").Value = 4.84
$ws.Range("AH5").Value = 10.27

# Remove the last data row (row 6) entirely
$ws.Rows("6:6").Delete()

# Tighten a handful of column widths (auto-fit side effect of shorter values)
$ws.Columns.Item(2).ColumnWidth = 6.17
$ws.Columns.Item(9).ColumnWidth = 6.17
$ws.Columns.Item(12).ColumnWidth = 6.17
$ws.Columns.Item(13).ColumnWidth = 6.17
$ws.Columns.Item(16).ColumnWidth = 6.17
$ws.Columns.Item(20).ColumnWidth = 7.17
$ws.Columns.Item(23).ColumnWidth = 6.17
$ws.Columns.Item(26).ColumnWidth = 6.17
$ws.Columns.Item(30).ColumnWidth = 6.17
$ws.Columns.Item(34).ColumnWidth = 6.17
